$d = $word.ActiveDocument

# Heading 3: apply Arial run font (rPr) on the style
$heading3 = $d.Styles("Heading3")
$heading3.Font.Name = "Arial"

# Indent 1 (custom style): apply Arial run font (rPr) on the style
$indent1 = $d.Styles("Indent1")
$indent1.Font.Name = "Arial"

# Subtitle: following paragraph should default back to Normal style
$subtitle = $d.Styles("Subtitle")
$subtitle.NextParagraphStyle = "Normal"

# Recitals (custom style): apply Arial run font (rPr) on the style
$recitals = $d.Styles("Recitals")
$recitals.Font.Name = "Arial"
